$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()

$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()

$ws.Range("H94").Value = 5260.364
$ws.Range("I94").Value = 4318.222
$ws.Range("K94").Value = 4318.222
$ws.Range("M94").Value = -3867.222

$ws.Range("H100").Value = 3688.9
$ws.Range("I100").Value = 3055.5715
$ws.Range("J100").Value = 5166.6665
$ws.Range("K100").Value = 3055.5715
$ws.Range("L100").Value = 5166.6665
$ws.Range("M100").Value = -2514.5715
$ws.Range("N100").Value = -6248.6665

$ws.Range("H116").Value = 3904.8
$ws.Range("I116").Value = 3881
$ws.Range("K116").Value = 3881
$ws.Range("M116").Value = -439

$ws.Range("H138").Value = 4535.45
$ws.Range("J138").Value = 6730.385
$ws.Range("L138").Value = 20191.155
$ws.Range("N138").Value = -30471.155

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4651.2
$ws.Range("I2").Value = 3310.75
$ws.Range("K2").Value = 3310.75
$ws.Range("M2").Value = -3197.75

$ws.Range("H44").Value = 11352.857
$ws.Range("J44").Value = 11620.5
$ws.Range("L44").Value = 11620.5
$ws.Range("N44").Value = -12596.5

$ws.Range("H97").Value = 825
$ws.Range("I97").Value = 810.3077
$ws.Range("K97").Value = 810.3077
$ws.Range("M97").Value = -314.3077

$ws.Range("H116").Value = 4651.2
$ws.Range("I116").Value = 3310.75
$ws.Range("K116").Value = 3310.75
$ws.Range("M116").Value = -1016.75

$ws.Range("H122").Value = 1249.75
$ws.Range("I122").Value = 1249.75
$ws.Range("K122").Value = 3749.25
$ws.Range("M122").Value = -1299.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4651.2
$ws.Range("I3").Value = 3310.75
$ws.Range("K3").Value = 3310.75
$ws.Range("M3").Value = -3196.75

$ws.Range("H86").Value = 6756.4443
$ws.Range("I86").Value = 3267
$ws.Range("K86").Value = 3267
$ws.Range("M86").Value = -2144

$ws.Range("H89").Value = 6756.4443
$ws.Range("I89").Value = 3267
$ws.Range("K89").Value = 16335
$ws.Range("M89").Value = -10719

$ws.Range("H94").Value = 472.4
$ws.Range("I94").Value = 472.4
$ws.Range("K94").Value = 472.4
$ws.Range("M94").Value = -21.39999999999998

$ws.Range("H99").Value = 771.0476
$ws.Range("I99").Value = 754.05
$ws.Range("K99").Value = 754.05
$ws.Range("M99").Value = 743.95

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 383.33334
$ws.Range("J6").Value = 200
$ws.Range("L6").Value = 200
$ws.Range("N6").Value = -426

$ws.Range("H22").Value = 1747.3334
$ws.Range("I22").Value = 399
$ws.Range("K22").Value = 399
$ws.Range("M22").Value = -49

$ws.Range("H23").Value = 9000
$ws.Range("I23").Value = 9000
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 9000
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -8760
$ws.Range("N23").ClearContents()

$ws.Range("H27").Value = 9000
$ws.Range("I27").Value = 9000
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 9000
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -8808
$ws.Range("N27").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2326.7273
$ws.Range("J5").Value = 2798.5
$ws.Range("L5").Value = 8395.5
$ws.Range("N5").Value = -8619.5

$ws.Range("H12").Value = 179.83333
$ws.Range("J12").Value = 214.88889
$ws.Range("L12").Value = 644.6666700000001
$ws.Range("N12").Value = -990.6666700000001

$ws.Range("H34").Value = 976.8570999999999
$ws.Range("I34").Value = 154
$ws.Range("J34").Value = 1114
$ws.Range("K34").Value = 462
$ws.Range("L34").Value = 3342
$ws.Range("M34").Value = -378
$ws.Range("N34").Value = -3510

$ws.Range("H39").Value = 6143
$ws.Range("J39").Value = 6143
$ws.Range("L39").Value = 18429
$ws.Range("N39").Value = -19017

$ws.Range("H55").Value = 2046
$ws.Range("J55").Value = 2553.2307
$ws.Range("L55").Value = 7659.6921
$ws.Range("N55").Value = -8013.6921

$ws.Range("H135").Value = 2326.7273
$ws.Range("J135").Value = 2798.5
$ws.Range("L135").Value = 25186.5
$ws.Range("N135").Value = -30256.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2041.5
$ws.Range("I102").Value = 1270.7142
$ws.Range("K102").Value = 1270.7142
$ws.Range("M102").Value = 351.2858000000001

$ws.Range("H122").Value = 2725.75
$ws.Range("I122").Value = 2481.4
$ws.Range("J122").Value = 3133
$ws.Range("K122").Value = 7444.200000000001
$ws.Range("L122").Value = 9399
$ws.Range("M122").Value = -4994.200000000001
$ws.Range("N122").Value = -14299

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 8875
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 8875
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 8875
$ws.Range("N46").Value = -9251
$ws.Range("M46").ClearContents()

$ws.Range("H93").Value = 1322.0834
$ws.Range("I93").Value = 1322.0834
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1322.0834
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -74.08339999999998
$ws.Range("N93").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2600.6667
$ws.Range("I100").Value = 2600.6667
$ws.Range("K100").Value = 5201.3334
$ws.Range("M100").Value = -4660.3334

$ws.Range("H107").Value = 836.25
$ws.Range("I107").Value = 750
$ws.Range("K107").Value = 2250
$ws.Range("M107").Value = -330

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws.Range("H122").Value = 5095
$ws.Range("I122").Value = 4150
$ws.Range("K122").Value = 12450
$ws.Range("M122").Value = -10000
